# Update imputed performance table values per author's desktop edit.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2.09 (1.92 to 2.27)";      New = "2.09 (1.93 to 2.26)" },
    @{ Old = "1.35% (1.16% to 1.55%)";   New = "1.35% (1.16% to 1.54%)" },
    @{ Old = "0.63 (0.51 to 0.75)";      New = "0.6 (0.47 to 0.73)" },
    @{ Old = "-0.26 (-0.33 to -0.19)";   New = "-0.27 (-0.34 to -0.21)" },
    @{ Old = "3.94%";                    New = "3.96%" },
    @{ Old = "1.14 (1.06 to 1.22)";      New = "1.13 (1.06 to 1.21)" },
    @{ Old = "0.54% (0.24% to 0.85%)";   New = "0.53% (0.25% to 0.8%)" },
    @{ Old = "-0.07 (-0.17 to 0.04)";    New = "-0.09 (-0.2 to 0.02)" },
    @{ Old = "-0.28 (-0.33 to -0.23)";   New = "-0.29 (-0.34 to -0.24)" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2)
}
